$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Rebuild the "Project logo" paragraph (paragraph 2): its paragraph-mark
#    run properties lose <w:color>/<w:lang> and gain <w:noProof/>, and the
#    embedded picture's <a:blip> gains cstate="print". Everything else in
#    the paragraph (text runs, bookmark, drawing metrics, relationship id)
#    stays the same.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$r2.Collapse(0)

$para2Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' + `
  '<w:pPr><w:rPr><w:noProof/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Project' + [char]0x00A0 + '</w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>logo :</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="logo"/>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:noProof/></w:rPr><w:drawing>' + `
    '<wp:inline distT="0" distB="0" distL="0" distR="0">' + `
      '<wp:extent cx="266700" cy="285750"/>' + `
      '<wp:effectExtent l="19050" t="0" r="0" b="0"/>' + `
      '<wp:docPr id="2" name="Image 0" descr="template.png"/>' + `
      '<wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr>' + `
      '<a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' + `
        '<pic:pic><pic:nvPicPr><pic:cNvPr id="0" name="template.png"/><pic:cNvPicPr/></pic:nvPicPr>' + `
        '<pic:blipFill><a:blip r:embed="rId5" cstate="print"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill>' + `
        '<pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="266700" cy="285750"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr>' + `
        '</pic:pic></a:graphicData></a:graphic>' + `
    '</wp:inline>' + `
  '</w:drawing></w:r>' + `
  '<w:bookmarkEnd w:id="0"/>' + `
'</w:p>'

[void]$r2.InsertXML($para2Xml)

# ---------------------------------------------------------------------------
# 2) Append a brand-new paragraph after it: "Logo (orginal size) : " plus a
#    second copy of the same picture (original/full size), wrapped in its
#    own bookmark "originalSizeLogo". InsertXML replaces the content of the
#    (collapsed) range's own paragraph rather than appending after it, so a
#    fresh empty paragraph is created first via InsertParagraphAfter and a
#    throw-away trailing paragraph anchors it away from the very end of the
#    story while its content gets written, then that trailing paragraph is
#    removed again.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs(3)
$p3.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs(3)
$r3 = $p3.Range

$para3Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' + `
  '<w:pPr><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve">Logo (orginal size) : </w:t></w:r>' + `
  '<w:bookmarkStart w:id="1" w:name="originalSizeLogo"/>' + `
  '<w:r><w:rPr><w:noProof/></w:rPr><w:drawing>' + `
    '<wp:inline distT="0" distB="0" distL="0" distR="0">' + `
      '<wp:extent cx="266700" cy="285750"/>' + `
      '<wp:effectExtent l="19050" t="0" r="0" b="0"/>' + `
      '<wp:docPr id="1" name="Image 0" descr="template.png"/>' + `
      '<wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr>' + `
      '<a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' + `
        '<pic:pic><pic:nvPicPr><pic:cNvPr id="0" name="template.png"/><pic:cNvPicPr/></pic:nvPicPr>' + `
        '<pic:blipFill><a:blip r:embed="rId5" cstate="print"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill>' + `
        '<pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="266700" cy="285750"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr>' + `
        '</pic:pic></a:graphicData></a:graphic>' + `
    '</wp:inline>' + `
  '</w:drawing></w:r>' + `
  '<w:bookmarkEnd w:id="1"/>' + `
'</w:p>'

[void]$r3.InsertXML($para3Xml)

$p4 = $d.Paragraphs(4)
$p4.Range.Delete()

Write-Output ("paragraphs=" + $d.Paragraphs.Count)
